$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Candidate ID 231027165 -> 231102298)
$ws.Range("A2").Value = "OPVVc194"
$ws.Range("B2").Value = 231102298
$ws.Range("C2").Value = "vpngbgw97"
$ws.Range("D2").Value = "eW!83qD&"
$ws.Range("F2").Value = "NCzpKaYz"
$ws.Range("G2").Value = "sxoI"

# Row 3 (Candidate ID 231027164 -> 231102297)
$ws.Range("A3").Value = "AFObi434"
$ws.Range("B3").Value = 231102297
$ws.Range("C3").Value = "hlzakqg57"
$ws.Range("D3").Value = "t$%E8fW5"
$ws.Range("F3").Value = "FmSVFzMG"
$ws.Range("G3").Value = "dXvA"
